$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D column (structure) cells for existing rows 2-9 were missing the
# shared white-fill formatting that every other data column already has.
# Bring them into line by copying the format from an already-correct
# sibling cell (C2) instead of assigning a brand-new fill.
$ws.Range("C2").Copy()
$ws.Range("D2:D9").PasteSpecial(-4122)

# Add the new ERV-Delta.9 (Molossus molossus) record as row 10, matching
# the look of the existing rows by first copying row 9's formatting down.
$ws.Range("A9:K9").Copy()
$ws.Range("A10:K10").PasteSpecial(-4122)

$ws.Range("A10").Value = "ERV-Delta.9-MolMol"
$ws.Range("B10").Value = "ERV-Delta.9-MolMol"
$ws.Range("C10").Value = "reference"
$ws.Range("D10").Value = "provirus"
$ws.Range("E10").Value = "ERV-Delta.9-Molossus_molossus"
$ws.Range("F10").Value = "Orthoretrovirinae"
$ws.Range("G10").Value = "Clade II"
$ws.Range("H10").Value = "Deltaretrovirus"
$ws.Range("I10").Value = "Unclassified"
$ws.Range("J10").Value = "Molossus molossus"
$ws.Range("K10").Value = "Species"

# The host-group-name cell was typed without picking up the row's
# background formatting, so clear it back to the default style.
$ws.Range("J10").ClearFormats()

# Finally, leave the same cell selected as in the authored workbook.
$ws.Range("D15").Select()
